$d = $word.ActiveDocument
Write-Host $d.CustomXMLParts.Count
